$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icosl"
$ws.Range("C2").Value = "Ctla4"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9314806666666667
$ws.Range("H2").Value = 2.794442
$ws.Range("I2").Value = 0.06288211749152639
$ws.Range("J2").Value = 0.06288211749152638
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005070666666666667
$ws.Range("N2").Value = 0.015212
$ws.Range("O2").Value = 0.008865695667453655
$ws.Range("P2").Value = 0.008865695667453653
$ws.Range("Q2").Value = 0.004723227967111112
$ws.Range("R2").Value = 0.042509051704
$ws.Range("S2").Value = 0.0005574937166049372
$ws.Range("T2").Value = 0.000557493716604937

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icosl"
$ws.Range("C3").Value = "Ctla4"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9314806666666667
$ws.Range("H3").Value = 2.794442
$ws.Range("I3").Value = 0.06288211749152639
$ws.Range("J3").Value = 0.06288211749152638
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5668716666666667
$ws.Range("N3").Value = 1.700615
$ws.Range("O3").Value = 0.9911343043325463
$ws.Range("P3").Value = 0.9911343043325463
$ws.Range("Q3").Value = 0.5280299979811112
$ws.Range("R3").Value = 4.75226998183
$ws.Range("S3").Value = 0.06232462377492146
$ws.Range("T3").Value = 0.06232462377492144

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icosl"
$ws.Range("C4").Value = "Ctla4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.399281333333333
$ws.Range("H4").Value = 7.197844
$ws.Range("I4").Value = 0.1619699646990985
$ws.Range("J4").Value = 0.1619699646990985
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005070666666666667
$ws.Range("N4").Value = 0.015212
$ws.Range("O4").Value = 0.008865695667453655
$ws.Range("P4").Value = 0.008865695667453653
$ws.Range("Q4").Value = 0.01216595588088889
$ws.Range("R4").Value = 0.109493602928
$ws.Range("S4").Value = 0.001435976414290419
$ws.Range("T4").Value = 0.001435976414290419

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icosl"
$ws.Range("C5").Value = "Ctla4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.399281333333333
$ws.Range("H5").Value = 7.197844
$ws.Range("I5").Value = 0.1619699646990985
$ws.Range("J5").Value = 0.1619699646990985
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5668716666666667
$ws.Range("N5").Value = 1.700615
$ws.Range("O5").Value = 0.9911343043325463
$ws.Range("P5").Value = 0.9911343043325463
$ws.Range("Q5").Value = 1.360084608228889
$ws.Range("R5").Value = 12.24076147406
$ws.Range("S5").Value = 0.1605339882848081
$ws.Range("T5").Value = 0.1605339882848081

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Icosl"
$ws.Range("C6").Value = "Ctla4"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.416663
$ws.Range("H6").Value = 1.249989
$ws.Range("I6").Value = 0.02812796084553394
$ws.Range("J6").Value = 0.02812796084553394
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005070666666666667
$ws.Range("N6").Value = 0.015212
$ws.Range("O6").Value = 0.008865695667453655
$ws.Range("P6").Value = 0.008865695667453653
$ws.Range("Q6").Value = 0.002112759185333333
$ws.Range("R6").Value = 0.019014832668
$ws.Range("S6").Value = 0.0002493739406025563
$ws.Range("T6").Value = 0.0002493739406025562

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Icosl"
$ws.Range("C7").Value = "Ctla4"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.416663
$ws.Range("H7").Value = 1.249989
$ws.Range("I7").Value = 0.02812796084553394
$ws.Range("J7").Value = 0.02812796084553394
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5668716666666667
$ws.Range("N7").Value = 1.700615
$ws.Range("O7").Value = 0.9911343043325463
$ws.Range("P7").Value = 0.9911343043325463
$ws.Range("Q7").Value = 0.2361944492483333
$ws.Range("R7").Value = 2.125750043235
$ws.Range("S7").Value = 0.02787858690493139
$ws.Range("T7").Value = 0.02787858690493138

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Icosl"
$ws.Range("C8").Value = "Ctla4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.0657
$ws.Range("H8").Value = 33.1971
$ws.Range("I8").Value = 0.7470199569638412
$ws.Range("J8").Value = 0.7470199569638412
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.005070666666666667
$ws.Range("N8").Value = 0.015212
$ws.Range("O8").Value = 0.008865695667453655
$ws.Range("P8").Value = 0.008865695667453653
$ws.Range("Q8").Value = 0.05611047613333334
$ws.Range("R8").Value = 0.5049942852
$ws.Range("S8").Value = 0.006622851595955742
$ws.Range("T8").Value = 0.006622851595955741

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Icosl"
$ws.Range("C9").Value = "Ctla4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.0657
$ws.Range("H9").Value = 33.1971
$ws.Range("I9").Value = 0.7470199569638412
$ws.Range("J9").Value = 0.7470199569638412
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5668716666666667
$ws.Range("N9").Value = 1.700615
$ws.Range("O9").Value = 0.9911343043325463
$ws.Range("P9").Value = 0.9911343043325463
$ws.Range("Q9").Value = 6.272831801833333
$ws.Range("R9").Value = 56.4554862165
$ws.Range("S9").Value = 0.7403971053678854
$ws.Range("T9").Value = 0.7403971053678854

